$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'Paris, David C. “Moral Education and the ‘Tie That Binds’ in Liberal Political Theory.” American Political Science Review 85, no. 3 (1991): 875–901. https://doi.org/10.2307/1963854.'
$ws.Cells.Item(2,2).Value = 'I examine the search for a “tie that binds, “ or “core” values, in liberal political theory, specifically Rawls''s recent arguments, and in proposals concerning moral education in the public schools. Both Rawls and the proponents of moral education appeal to consensus or shared values, but the search for core values in both theory and practice is only partly successful. Specifically, this search is misguided insofar as it does not reflect how values are embedded in specific institutions and practices. The various forms of moral education in the public schools, both implicit and explicit, illustrate a consensus about a range of moral and intellectual virtues that is broader and more complex than arguments for core values allow. Comparing arguments concerning core values in political theory and moral education suggests how liberal political theory might deal with questions of consensus, justification, and the task of political theory generally.'
$ws.Cells.Item(2,3).Value = 'Hamilton College'
$ws.Cells.Item(2,4).Value = ''

# Row 3
$ws.Cells.Item(3,1).Value = 'Richardson, Bradley M. “European Party Loyalties Revisited.” American Political Science Review 85, no. 3 (1991): 751–75. https://doi.org/10.2307/1963849.'
$ws.Cells.Item(3,2).Value = 'The likelihood that party loyalties reflect the characteristics of parties and political competition within specific political environments was generally ignored in critiques of party identification abroad. In this article, I identify substantial frequencies of stable partisan ties among British, Dutch, and German voters in the 1960s through the 1980s. Stable party loyalists in turn usually vote for their preferred party in subsequent elections. Long-term partisan ties are most common among supporters of traditional cleavage parties. These long-term party loyalties reflect the effects of both long-standing hostilities toward opposing parties and well-developed party principles. Since partisan responses are internally consistent and subsume ideational as well as affective components, these feelings resemble social psychology''s affect-laden schemata.'
$ws.Cells.Item(3,3).Value = 'Ohio State University'
$ws.Cells.Item(3,4).Value = ''

# Row 4
$ws.Cells.Item(4,1).Value = 'Wagner, R. Harrison. “Nuclear Deterrence, Counterforce Strategies, and the Incentive to Strike First.” American Political Science Review 85, no. 3 (1991): 727–49. https://doi.org/10.2307/1963848.'
$ws.Cells.Item(4,2).Value = 'I use the theory of games to investigate issues about how to understand the use of nuclear counterforce strategies by the United States and the Soviet Union during the Cold War. The function of the counterforce strategies I model is not to enable a state confidently to launch a nuclear attack but to convince its adversary that the probability that it might do so as a last resort is greater than zero. The models allow one to investigate rational behavior when information is incomplete and there is an incentive to strike first, and therefore provide a way to explore controversies about the effect of counterforce strategies on both the credibility of extended deterrence and the possibility of inadvertent nuclear war. The models suggest, contrary to the claims of a number of writers, that the use of nuclear counterforce strategies is not necessarily inconsistent with rational behavior and provide some insight into the relation between counterforce strategies and brinkmanship models of deterrence.'
$ws.Cells.Item(4,3).Value = 'University of Texas Austin'
$ws.Cells.Item(4,4).Value = ''

# Row 5
$ws.Cells.Item(5,1).Value = 'Quinn, Dennis P., and Robert Y. Shapiro. “Business Political Power: The Case of Taxation.” American Political Science Review 85, no. 3 (1991): 851–74. https://doi.org/10.2307/1963853.'
$ws.Cells.Item(5,2).Value = 'We examine contending views about the forms and mechanisms of business power in U.S. politics by estimating time series models explaining taxation and redistribution. Taxation and redistribution constitute strong cases for theories about business and class power, since all firms have an interest in reducing taxation. We find that changes in corporate taxation and in redistribution between capital gains income and earned income and between corporate taxation and individual taxation are strongly influenced by political partisanship, with Democratic administrations increasing the tax burden on firms and their owners. How far corporations engage in electoral financing—measured through the establishment of corporate political action committees—is also influential. The models show some evidence consistent with the understanding of class power conceptualized by Bowles, Gordon, and Weisskopf and the presence of election cycle effects but inconsistent with implications arising from the structural dependence of the state on capital and asset concentration in the largest corporations as a mechanism of class power.'
$ws.Cells.Item(5,3).Value = 'Georgetown University'
$ws.Cells.Item(5,4).Value = 'Columbia University'

# Row 6
$ws.Cells.Item(6,1).Value = 'Scholz, John T., Jim Twombly, and Barbara Headrick. “Street-Level Political Controls Over Federal Bureaucracy.” American Political Science Review 85, no. 3 (1991): 829–50. https://doi.org/10.2307/1963852.'
$ws.Cells.Item(6,2).Value = 'Local partisan activities of legislators and their electoral coalitions systematically influence field office activities of federal bureaucracies in their electoral districts. This alternative to centralized democratic controls over bureaucracy occurs because discretionary policy decisions made at the field office level are influenced by local resources generated through partisan activities. Our study of county-level Occupational Safety and Health Administration enforcement in New York (1976–85) finds that county, state, and federal elected officials influence local enforcement activities, with liberal, Democratic legislators associated with more active enforcement. The county political parties are most influential for activities with the most local discretion, while members of Congress are more influential for local activities more readily controlled by the national office.'
$ws.Cells.Item(6,3).Value = 'State University of New York Stony Brook'
$ws.Cells.Item(6,4).Value = 'State University of New York Buffalo / Texas A&M University'

# Row 7
$ws.Cells.Item(7,1).Value = 'Snidal, Duncan. “Relative Gains and the Pattern of International Cooperation.” American Political Science Review 85, no. 3 (1991): 701–26. https://doi.org/10.2307/1963847.'
$ws.Cells.Item(7,2).Value = 'Many political situations involve competitions where winning is more important than doing well. In international politics, this relative gains problem is widely argued to be a major impediment to cooperation under anarchy. After discussing why states might seek relative gains, I demonstrate that the hypothesis holds very different implications from those usually presumed. Relative gains do impede cooperation in the two-actor case and provide an important justification for treating international anarchy as a prisoner''s dilemma problem; but if the initial absolute gains situation is not a prisoner''s dilemma, relative gains seeking is much less consequential. Its significance is even more attenuated with more than two competitors. Relative gains cannot prop up the realist critique of international cooperation theory, but may affect the pattern of cooperation when a small number of states are the most central international actors.'
$ws.Cells.Item(7,3).Value = 'University of Chicago'
$ws.Cells.Item(7,4).Value = ''

# Row 8
$ws.Cells.Item(8,1).Value = 'Remmer, Karen L. “The Political Impact of Economic Crisis in Latin America in the 1980s.” American Political Science Review 85, no. 3 (1991): 777–800. https://doi.org/10.2307/1963850.'
$ws.Cells.Item(8,2).Value = 'Research on the political implications of economic conditions is separated into two relatively distinct bodies of literature. I bridge the theoretical gap between them by examining the effects of economic crisis on electoral outcomes in Latin America from 1982 to 1990. An analysis of 21 competitive elections indicates that crisis conditions undermine support for incumbents and provoke high levels of electoral volatility but without necessarily fostering the growth of political extremism or the exhaustion of elite consensus associated with the breakdown of democracy. The results also suggest that the relationship between economic conditions and electoral instability is mediated by party system structure rather than democratic age. Paradoxically, the findings buttress prior research on electoral outcomes in the comparatively stable and homogeneous Organization for Economic Cooperation and Development (OECD) nations while undercutting theoretical frameworks elaborated with specific reference to the breakdown and consolidation of Third World democracy.'
$ws.Cells.Item(8,3).Value = 'University of New Mexico'
$ws.Cells.Item(8,4).Value = ''

# Row 9
$ws.Cells.Item(9,1).Value = 'Wood, B. Dan, and Richard W. Waterman. “The Dynamics of Political Control of the Bureaucracy.” American Political Science Review 85, no. 3 (1991): 801–28. https://doi.org/10.2307/1963851.'
$ws.Cells.Item(9,2).Value = 'A new paradigm of political-bureaucratic relations emerged through the 1980s holding that U.S. democratic institutions continuously shape nonelective public bureaucracies. Several empirical studies support the paradigm with evidence suggestive of political manipulation but none reveals the scope or specific mechanisms of political control. We explore the dynamics of political control of the bureaucracy explicitly to determine the scope and mechanisms. We examine output time series from seven different public bureaucracies for responsiveness to political tools applied in the late Carter and early Reagan administrations. We find responsiveness in all seven cases. The evidence also shows that political appointments—a shared power of the president and Congress—is the most important instrument of political control; changing budgets, legislation, congressional signals, and administrative reorganizations are less important. These findings confirm intuitive assertions by institutional scholars and suggest a method of “policy monitoring” that could enhance future democratic control of the bureaucracy.'
$ws.Cells.Item(9,3).Value = 'Texas A&M University'
$ws.Cells.Item(9,4).Value = 'University of New Mexico'
